$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row (root / Todos / Todos las db y sus tablas.) after the header row ---
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "root"
$ws.Range("B2").Value = "Todos"
$ws.Range("C2").Value = "Todos las db y sus tablas."

# --- Row 4 (previously row 3): "Usuarios secundarios tipo A" gets new privileges/access/notes ---
$ws.Range("B4").Value = "Solo puede actualizar, insertar y consultar datos en las tablas ya creadas. (Insert, update y select)"
$ws.Range("C4").Value = "Todas la tablas pero no a las columnas con datos confidenciales del paciente."
$ws.Range("D4").Value = "La tabla ""pacientes"" contiene columnas confidenciales."
$ws.Rows.Item(4).RowHeight = 45

# Make part of B4 bold to reproduce the rich-text run split: "...creadas." + " (" (bold) + "Insert, update y select)"
$boldPart = $ws.Range("B4").Characters(76, 2)
$boldPart.Font.Bold = $true
$tailPart = $ws.Range("B4").Characters(78, 24)
$tailPart.Font.Bold = $false

# --- Row 5 (previously row 4): "Usuarios secundarios tipo B" gets new access/notes ---
$ws.Range("C5").Value = "Todas la tablas pero no a las columnas con datos confidenciales del paciente."
$ws.Range("D5").Value = "La tabla ""pacientes"" contiene columnas confidenciales."
$ws.Rows.Item(5).RowHeight = 45

$ws.Range("E4").Select()
